# strategy test iteration 4 player 0 update strategy
#
# 1. "regrets" sheet: move the view/selection from M6 -> B4 (and stop being
#    the active tab).
# 2. "action counts" sheet: the iteration-4 / player-0 strategy update fills
#    in previously-blank action-count cells with 0 (and the actual action
#    taken, fold on the river for sequence id 0, with a count of 1 in K4).
#    Also move the view/selection from B4 -> K5.
# 3. "sequence table" sheet: move the view/selection from A6 -> D4, and make
#    this the active tab/sheet (activeTab=2 in the saved workbook).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "regrets" sheet
# ---------------------------------------------------------------------------
$wsRegrets = $wb.Worksheets.Item("regrets")
$wsRegrets.Activate()
$wsRegrets.Range("B4").Select()

# ---------------------------------------------------------------------------
# "action counts" sheet
# ---------------------------------------------------------------------------
$wsActionCounts = $wb.Worksheets.Item("action counts")
$wsActionCounts.Activate()

$zeroCells = @(
    "B4", "C4", "D4", "E4", "F4", "G4", "H4", "I4", "J4", "L4",
    "M4", "B5", "C5", "D5", "E5", "F5", "G5", "H5", "I5", "J5",
    "K5", "L5", "M5", "B6", "C6", "E6", "F6", "H6", "I6", "K6",
    "L6", "C7", "D7", "F7", "G7", "I7", "J7", "L7", "M7", "B8",
    "C8", "E8", "F8", "H8", "I8", "K8", "L8", "B9", "C9", "E9",
    "F9", "H9", "I9", "K9", "L9", "B10", "C10", "E10", "F10", "H10",
    "I10", "K10", "L10", "B11", "C11", "E11", "F11", "H11", "I11", "K11",
    "L11", "B12", "C12", "D12", "E12", "F12", "G12", "H12", "I12", "J12",
    "K12", "L12", "M12", "C13", "D13", "F13", "G13", "I13", "J13", "L13",
    "M13", "B14", "C14", "E14", "F14", "H14", "I14", "K14", "L14", "B15",
    "C15", "E15", "F15", "H15", "I15", "K15", "L15", "B16", "C16", "E16",
    "F16", "H16", "I16", "K16", "L16", "B17", "C17", "E17", "F17", "H17",
    "I17", "K17", "L17", "C18", "D18", "F18", "G18", "I18", "J18", "L18",
    "M18", "B19", "C19", "E19", "F19", "H19", "I19", "K19", "L19", "B20",
    "C20", "E20", "F20", "H20", "I20", "K20", "L20", "B21", "C21", "E21",
    "F21", "H21", "I21", "K21", "L21"
)

foreach ($addr in $zeroCells) {
    $wsActionCounts.Range($addr).Value = 0
}

# K4 = fold count on the river (round index 3) for sequence id 0 -> player 0
# took that action once.
$wsActionCounts.Range("K4").Value = 1

$wsActionCounts.Range("K5").Select()

# ---------------------------------------------------------------------------
# "sequence table" sheet (left as the active sheet/tab after the edit)
# ---------------------------------------------------------------------------
$wsSequenceTable = $wb.Worksheets.Item("sequence table")
$wsSequenceTable.Activate()
$wsSequenceTable.Range("D4").Select()
